$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '68.060.47'
$ws.Cells.Item(2, 5).Value = '  -3.08%  '

$ws.Cells.Item(3, 4).Value = '3.843.24'
$ws.Cells.Item(3, 5).Value = '  -2.41%  '

$ws.Cells.Item(4, 5).Value = '  -0.10%  '

$ws.Cells.Item(5, 4).Value = '600.71'
$ws.Cells.Item(5, 5).Value = '  -1.73%  '

$ws.Cells.Item(6, 4).Value = '167.78'
$ws.Cells.Item(6, 5).Value = '  -2.30%  '

$ws.Cells.Item(7, 4).Value = '3.843.02'
$ws.Cells.Item(7, 5).Value = '  -2.39%  '

$ws.Cells.Item(8, 5).Value = '  -0.02%  '

$ws.Cells.Item(9, 4).Value = '0.529'
$ws.Cells.Item(9, 5).Value = '  -2.01%  '

$ws.Cells.Item(10, 5).Value = '  -4.40%  '

$ws.Cells.Item(11, 4).Value = '6.47'
$ws.Cells.Item(11, 5).Value = '  +0.04%  '

$ws.Cells.Item(12, 5).Value = '  -3.00%  '

$ws.Cells.Item(13, 5).Value = '  +0.67%  '

$ws.Cells.Item(14, 4).Value = '37.03'
$ws.Cells.Item(14, 5).Value = '  -4.36%  '

$ws.Cells.Item(15, 4).Value = '4.482.49'
$ws.Cells.Item(15, 5).Value = '  -2.55%  '

$ws.Cells.Item(16, 4).Value = '3.840.66'
$ws.Cells.Item(16, 5).Value = '  -3.24%  '

$ws.Cells.Item(17, 4).Value = '68.137.91'
$ws.Cells.Item(17, 5).Value = '  -2.97%  '

$ws.Cells.Item(18, 4).Value = '18.31'
$ws.Cells.Item(18, 5).Value = '  -1.75%  '

$ws.Cells.Item(19, 4).Value = '7.40'
$ws.Cells.Item(19, 5).Value = '  -3.95%  '

$ws.Cells.Item(20, 5).Value = '  -0.85%  '

$ws.Cells.Item(21, 4).Value = '11.06'
$ws.Cells.Item(21, 5).Value = '  -0.19%  '

$ws.Cells.Item(22, 4).Value = '466.55'
$ws.Cells.Item(22, 5).Value = '  -6.22%  '

$ws.Cells.Item(23, 4).Value = '0.735'
$ws.Cells.Item(23, 5).Value = '  -2.08%  '

$ws.Cells.Item(24, 5).Value = '  -4.05%  '

$ws.Cells.Item(25, 5).Value = '  -3.75%  '

$ws.Cells.Item(26, 5).Value = '  -3.13%  '

$ws.Cells.Item(27, 4).Value = '12.12'
$ws.Cells.Item(27, 5).Value = '  -2.48%  '

$ws.Cells.Item(28, 4).Value = '10.07'
$ws.Cells.Item(28, 5).Value = '  -1.51%  '

$ws.Cells.Item(29, 5).Value = '  -0.22%  '

$ws.Cells.Item(30, 5).Value = '  -1.86%  '

$ws.Cells.Item(31, 4).Value = '3.991.19'
$ws.Cells.Item(31, 5).Value = '  -2.43%  '

$ws.Cells.Item(32, 4).Value = '7.65'
$ws.Cells.Item(32, 5).Value = '  -3.37%  '

$ws.Cells.Item(33, 4).Value = '2.32'
$ws.Cells.Item(33, 5).Value = '  -6.08%  '

$ws.Cells.Item(34, 4).Value = '31.28'
$ws.Cells.Item(34, 5).Value = '  -3.53%  '

$ws.Cells.Item(35, 4).Value = '9.57'
$ws.Cells.Item(35, 5).Value = '  -0.96%  '

$ws.Cells.Item(36, 4).Value = '3.803.14'
$ws.Cells.Item(36, 5).Value = '  -2.55%  '

$ws.Cells.Item(37, 5).Value = '  -3.85%  '

$ws.Cells.Item(38, 4).Value = '3.64'
$ws.Cells.Item(38, 5).Value = '  +9.62%  '

$ws.Cells.Item(39, 5).Value = '  -1.01%  '

$ws.Cells.Item(40, 5).Value = '  -3.32%  '

$ws.Cells.Item(41, 4).Value = '5.92'
$ws.Cells.Item(41, 5).Value = '  -4.43%  '

$ws.Cells.Item(42, 5).Value = '  -0.11%  '

$ws.Cells.Item(43, 4).Value = '0.314'
$ws.Cells.Item(43, 5).Value = '  -5.12%  '

$ws.Cells.Item(44, 5).Value = '  -7.26%  '

$ws.Cells.Item(45, 4).Value = '421.84'
$ws.Cells.Item(45, 5).Value = '  -4.06%  '

$ws.Cells.Item(46, 5).Value = '  -0.13%  '

$ws.Cells.Item(47, 5).Value = '  -0.01%  '

$ws.Cells.Item(48, 4).Value = '0.000290'
$ws.Cells.Item(48, 5).Value = '  +4.41%  '

$ws.Cells.Item(49, 4).Value = '46.95'
$ws.Cells.Item(49, 5).Value = '  -2.95%  '

$ws.Cells.Item(50, 4).Value = '142.37'
$ws.Cells.Item(50, 5).Value = '  -0.80%  '

$ws.Cells.Item(51, 2).Value = 'EnergySwap'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(51, 4).Value = '26.21'
$ws.Cells.Item(51, 5).Value = '  +2.31%  '

Write-Output "Applied cryptos update"
